# Refresh the crypto price/volume snapshot (Price column D, Volume(1h)
# column E) with the latest scraped figures.
#
# Commit message: "Updated cryptos list on Sat Oct 28 20:40:55 UTC 2023
# with GitHub Actions"
#
# Every cell in D2:E51 is stored as plain text in this workbook (even the
# ones that look like plain decimals, e.g. "226.14"), so a handful of the
# new Price values would otherwise be auto-coerced to numbers by Excel's
# smart-entry parsing. For those we lead with an apostrophe (forces text
# entry, like a user typing '226.21 in the UI) and then reset the cell
# Style back to Normal so no stray number-format/quote-prefix styling is
# left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '34.170.30'
$ws.Range('E2').Value = '  +1.19%  '
$ws.Range('D3').Value = '1.781.50'
$ws.Range('E3').Value = '  +0.13%  '
$ws.Range('E4').Value = '  +0.28%  '
$ws.Range('D5').Value = "'226.21"
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.01%  '
$ws.Range('E6').Value = '  +0.19%  '
$ws.Range('E7').Value = '  +0.33%  '
$ws.Range('E8').Value = '  -1.23%  '
$ws.Range('E9').Value = '  +0.69%  '
$ws.Range('E10').Value = '  +1.33%  '
$ws.Range('E11').Value = '  +1.20%  '
$ws.Range('D12').Value = '2.037.82'
$ws.Range('E12').Value = '  +0.20%  '
$ws.Range('D13').Value = "'10.99"
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -1.77%  '
$ws.Range('D14').Value = '1.774.99'
$ws.Range('E14').Value = '  +0.06%  '
$ws.Range('D15').Value = '34.134.41'
$ws.Range('E16').Value = '  +1.81%  '
$ws.Range('E17').Value = '  +1.13%  '
$ws.Range('D18').Value = "'67.97"
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.04%  '
$ws.Range('D19').Value = "'247.13"
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +3.66%  '
$ws.Range('D20').Value = '0.0₃0793'
$ws.Range('E20').Value = '  +2.61%  '
$ws.Range('E21').Value = '  +3.69%  '
$ws.Range('E22').Value = '  +0.23%  '
$ws.Range('E23').Value = '  +2.10%  '
$ws.Range('D25').Value = "'162.74"
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +1.46%  '
$ws.Range('D26').Value = "'7.19"
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +2.27%  '
$ws.Range('D27').Value = "'16.31"
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.23%  '
$ws.Range('E28').Value = '  +1.52%  '
$ws.Range('E30').Value = '  +0.41%  '
$ws.Range('E31').Value = '  +1.54%  '
$ws.Range('D32').Value = "'3.74"
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +3.92%  '
$ws.Range('E33').Value = '  +4.88%  '
$ws.Range('E34').Value = '  -0.22%  '
$ws.Range('D35').Value = '1.442.15'
$ws.Range('D36').Value = "'0.654"
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +1.00%  '
$ws.Range('E37').Value = '  +7.10%  '
$ws.Range('E38').Value = '  +3.38%  '
$ws.Range('D39').Value = "'1.05"
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.42%  '
$ws.Range('D40').Value = "'80.17"
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +2.42%  '
$ws.Range('E41').Value = '  +0.77%  '
$ws.Range('E42').Value = '  +1.44%  '
$ws.Range('D43').Value = "'13.71"
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.71%  '
$ws.Range('E44').Value = '  +0.61%  '
$ws.Range('E45').Value = '  +2.11%  '
$ws.Range('E46').Value = '  +3.66%  '
$ws.Range('E48').Value = '  +0.45%  '
$ws.Range('D49').Value = '1.940.10'
$ws.Range('E49').Value = '  +0.69%  '
$ws.Range('D50').Value = "'104.62"
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.81%  '
$ws.Range('E51').Value = '  +0.26%  '
